$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -35.24428059585877
$ws.Range("C2").Value = 269.5602310681157

$ws.Range("B3").Value = -43.81402676005936
$ws.Range("C3").Value = -54.68498494151564

$ws.Range("B4").Value = -38.07047593582452
$ws.Range("C4").Value = -98.15290549268987

$ws.Range("B5").Value = -44.07929141415514
$ws.Range("C5").Value = 52.62759760343716

$ws.Range("B6").Value = -29.18970863480975
$ws.Range("C6").Value = -86.42868143611607

$ws.Range("B8").Value = -36.91842092174158
$ws.Range("C8").Value = 153.3741001090112

$ws.Range("B9").Value = -37.84140403772091
$ws.Range("C9").Value = -79.12443946814697

$ws.Range("B10").Value = -26.48295957781719
$ws.Range("C10").Value = -77.64233321010346
